$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update "想去人数" (column F) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 550
$ws1.Range("F4").Value = 196
$ws1.Range("F6").Value = 505
$ws1.Range("F7").Value = 104
$ws1.Range("F8").Value = 119
$ws1.Range("F9").Value = 44
$ws1.Range("F10").Value = 6732
$ws1.Range("F12").Value = 371
$ws1.Range("F13").Value = 3024
$ws1.Range("F14").Value = 195
$ws1.Range("F15").Value = 341
$ws1.Range("F17").Value = 542

# Sheet "全部类型" (fourth sheet) - update "想去人数" (column F) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 550
$ws4.Range("F6").Value = 196
$ws4.Range("F8").Value = 505
$ws4.Range("F9").Value = 104
$ws4.Range("F10").Value = 119
$ws4.Range("F11").Value = 44
$ws4.Range("F13").Value = 6732
$ws4.Range("F16").Value = 371
$ws4.Range("F17").Value = 3024
$ws4.Range("F18").Value = 195
$ws4.Range("F19").Value = 341
$ws4.Range("F21").Value = 542
